$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34-92 down to 35-93
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly record
$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 44775
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 100114007
$ws.Range("G34").Value = "Jengibre"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 11000
$ws.Range("L34").Value = 12000
$ws.Range("M34").Value = 11400
$ws.Range("N34").Value = "$/caja 13 kilos"
$ws.Range("O34").Value = "Perú"
$ws.Range("P34").Value = 877
$ws.Range("Q34").Value = 13
$ws.Range("R34").Value = "Hortaliza"
